# StructureDefinition-TiposVacunaRNI.xlsx — "limpiando y borrando lo maximo"
#
# 1) Fix casing of the CoreCL path segment ("corecl" -> "CoreCL") in the
#    canonical StructureDefinition URL and the ValueSet binding URL.
# 2) Refresh the generation Date stamp.
# 3) The above text edits change rendered string widths, so the IG
#    publisher's spreadsheet step re-autofit the "Elements" sheet's
#    columns; reproduce the resulting column widths as closely as this
#    runtime's ColumnWidth (whole-pixel, 1/6-character granularity)
#    allows.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------------
$ws1.Range("B2").Value = "https://hl7chile.cl/fhir/ig/CoreCL/StructureDefinition/TiposVacunaRNI"
$ws1.Range("B8").Value = "2022-12-12T20:08:16-03:00"

# --- Elements sheet --------------------------------------------------------
# Q5 ("Base Definition" Source row) shares the same StructureDefinition URL
# string as Metadata!B2 in the workbook's shared-string table, so it must be
# updated in lockstep to keep both cells pointing at identical text.
$ws2.Range("Q5").Value = "https://hl7chile.cl/fhir/ig/CoreCL/StructureDefinition/TiposVacunaRNI"
$ws2.Range("Y7").Value = "https://hl7chile.cl/fhir/ig/CoreCL/ValueSet/VSTiposVacunas"

# Recomputed column widths (character units) for columns A..AJ.
$widths = @(
    18.170572916666668,
    10.299479166666666,
    6.850260416666667,
    5.069010416666667,
    3.8658854166666665,
    4.240885416666667,
    13.791666666666666,
    11.154947916666666,
    11.912760416666666,
    19.869791666666668,
    35.576822916666664,
    99.86979166666667,
    99.86979166666667,
    99.86979166666667,
    12.694010416666666,
    19.869791666666668,
    19.869791666666668,
    19.869791666666668,
    19.869791666666668,
    14.881510416666666,
    15.256510416666666,
    16.244791666666668,
    15.479166666666666,
    18.080729166666668,
    53.647135416666664,
    4.858072916666667,
    18.897135416666668,
    39.205729166666664,
    14.154947916666666,
    11.471354166666666,
    16.893229166666668,
    8.666666666666666,
    9.045572916666666,
    11.385416666666666,
    99.86979166666667,
    21.842447916666668
)

for ($i = 0; $i -lt $widths.Length; $i++) {
    $ws2.Columns.Item($i + 1).ColumnWidth = $widths[$i]
}
